$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1635220125786163
$ws.Range("C2").Value = 0.6163522012578616
$ws.Range("J2").Value = 0.02515723270440252
$ws.Range("P2").Value = 0.1132075471698113
$ws.Range("S2").Value = 0.08176100628930817
$ws.Range("B3").Value = 0.02777777777777778
$ws.Range("C3").Value = 0.04629629629629629
$ws.Range("J3").Value = 0.03703703703703703
$ws.Range("S3").Value = 0.2222222222222222
$ws.Range("J4").Value = 0.04545454545454546
$ws.Range("S4").Value = 0.4090909090909091
$ws.Range("P5").Value = 0.5
$ws.Range("S5").Value = 0.5
$ws.Range("B6").Value = 0.07920792079207921
$ws.Range("D6").Value = 0.009900990099009901
$ws.Range("F6").Value = 0.09900990099009901
$ws.Range("J6").Value = 0.2376237623762376
$ws.Range("O6").Value = 0.0396039603960396
$ws.Range("Q6").Value = 0.1683168316831683
$ws.Range("R6").Value = 0.0594059405940594
$ws.Range("S6").Value = 0.3069306930693069
$ws.Range("B7").Value = 0.1012658227848101
$ws.Range("F7").Value = 0.06329113924050633
$ws.Range("J7").Value = 0.1139240506329114
$ws.Range("O7").Value = 0.0379746835443038
$ws.Range("Q7").Value = 0.1139240506329114
$ws.Range("R7").Value = 0.1518987341772152
$ws.Range("S7").Value = 0.4177215189873418
$ws.Range("B8").Value = 0.1288659793814433
$ws.Range("D8").Value = 0.03608247422680412
$ws.Range("E8").Value = 0.005154639175257732
$ws.Range("F8").Value = 0.05154639175257732
$ws.Range("J8").Value = 0.07731958762886598
$ws.Range("O8").Value = 0.02577319587628866
$ws.Range("Q8").Value = 0.1701030927835052
$ws.Range("R8").Value = 0.1237113402061856
$ws.Range("S8").Value = 0.3814432989690721
$ws.Range("B9").Value = 0.08888888888888889
$ws.Range("D9").Value = 0.01111111111111111
$ws.Range("F9").Value = 0.06666666666666667
$ws.Range("J9").Value = 0.1666666666666667
$ws.Range("O9").Value = 0.01111111111111111
$ws.Range("Q9").Value = 0.1888888888888889
$ws.Range("R9").Value = 0.05555555555555555
$ws.Range("S9").Value = 0.4111111111111111
$ws.Range("B10").Value = 0.1359867330016584
$ws.Range("D10").Value = 0.02321724709784411
$ws.Range("E10").Value = 0.001658374792703151
$ws.Range("F10").Value = 0.06633499170812604
$ws.Range("J10").Value = 0.1326699834162521
$ws.Range("O10").Value = 0.01658374792703151
$ws.Range("Q10").Value = 0.1923714759535655
$ws.Range("R10").Value = 0.07131011608623548
$ws.Range("S10").Value = 0.3598673300165838
$ws.Range("G11").Value = 0.1484375
$ws.Range("J11").Value = 0.09375
$ws.Range("K11").Value = 0.1953125
$ws.Range("L11").Value = 0.5390625
$ws.Range("S11").Value = 0.0234375
$ws.Range("G12").Value = 0.75
$ws.Range("J12").Value = 0.1527777777777778
$ws.Range("K12").Value = 0.02777777777777778
$ws.Range("L12").Value = 0.04166666666666666
$ws.Range("S12").Value = 0.02777777777777778
$ws.Range("G13").Value = 0.5
$ws.Range("J13").Value = 0.5
$ws.Range("F15").Value = 0.01941747572815534
$ws.Range("H15").Value = 0.1650485436893204
$ws.Range("I15").Value = 0.07766990291262135
$ws.Range("J15").Value = 0.3300970873786408
$ws.Range("K15").Value = 0.06796116504854369
$ws.Range("O15").Value = 0.04854368932038835
$ws.Range("S15").Value = 0.2912621359223301
$ws.Range("H16").Value = 0.1717171717171717
$ws.Range("I16").Value = 0.0505050505050505
$ws.Range("J16").Value = 0.3939393939393939
$ws.Range("K16").Value = 0.101010101010101
$ws.Range("M16").Value = 0.04040404040404041
$ws.Range("O16").Value = 0.0303030303030303
$ws.Range("S16").Value = 0.2121212121212121
$ws.Range("F17").Value = 0.02659574468085106
$ws.Range("H17").Value = 0.2180851063829787
$ws.Range("I17").Value = 0.101063829787234
$ws.Range("J17").Value = 0.398936170212766
$ws.Range("K17").Value = 0.06914893617021277
$ws.Range("M17").Value = 0.02659574468085106
$ws.Range("O17").Value = 0.05851063829787234
$ws.Range("S17").Value = 0.101063829787234
$ws.Range("F18").Value = 0.02197802197802198
$ws.Range("H18").Value = 0.1428571428571428
$ws.Range("I18").Value = 0.1208791208791209
$ws.Range("J18").Value = 0.3736263736263736
$ws.Range("K18").Value = 0.04395604395604396
$ws.Range("M18").Value = 0.01098901098901099
$ws.Range("O18").Value = 0.1098901098901099
$ws.Range("S18").Value = 0.1758241758241758
$ws.Range("F19").Value = 0.01706484641638225
$ws.Range("H19").Value = 0.1843003412969283
$ws.Range("I19").Value = 0.08020477815699659
$ws.Range("J19").Value = 0.4112627986348123
$ws.Range("K19").Value = 0.1143344709897611
$ws.Range("M19").Value = 0.0136518771331058
$ws.Range("N19").Value = 0.001706484641638225
$ws.Range("O19").Value = 0.06484641638225255
$ws.Range("S19").Value = 0.1126279863481229

Write-Output "Applied 105 cell updates"
